$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.103.54'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '3.792.03'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '599.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.33'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.71%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.516'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.68%  '
$ws.Range('E9').Value = '  -1.24%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.54'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.67%  '
$ws.Range('E12').Value = '  -2.24%  '
$ws.Range('E13').Value = '  -0.83%  '
$ws.Range('D14').Value = '4.427.22'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').Value = '3.787.09'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').Value = '68.094.08'
$ws.Range('E16').Value = '  +0.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.30'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('E18').Value = '  +2.23%  '
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.47'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('E21').Value = '  -2.56%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000145'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.55%  '
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.96'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.14%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').Value = '3.939.61'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('E31').Value = '  -4.95%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.21'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0996'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.28'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.90%  '
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.985'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.42%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  -1.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.28'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.297'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.92%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '152.49'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.03%  '
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.71'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.34%  '
$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '388.85'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.52%  '
